$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 4573.75
$ws.Range("I33").Value = 5134.125
$ws.Range("K33").Value = 5134.125
$ws.Range("M33").Value = -4905.125
$ws.Range("H40").Value = 5288.2
$ws.Range("I40").Value = 3583
$ws.Range("K40").Value = 3583
$ws.Range("M40").Value = -3408
$ws.Range("H41").Value = 1334.3846
$ws.Range("I41").Value = 1050
$ws.Range("K41").Value = 1050
$ws.Range("M41").Value = -610
$ws.Range("H51").Value = 4849.7144
$ws.Range("I51").Value = 5399.5
$ws.Range("K51").Value = 5399.5
$ws.Range("M51").Value = -4915.5
$ws.Range("H74").Value = 4851.909
$ws.Range("I74").Value = 4538.6
$ws.Range("K74").Value = 4538.6
$ws.Range("M74").Value = -3602.6
$ws.Range("H77").Value = 4851.909
$ws.Range("I77").Value = 4538.6
$ws.Range("K77").Value = 22693
$ws.Range("M77").Value = -18013
$ws.Range("H103").Value = 621.7692
$ws.Range("I103").Value = 759.3333
$ws.Range("J103").Value = 312.25
$ws.Range("K103").Value = 2277.9999
$ws.Range("L103").Value = 936.75
$ws.Range("M103").Value = -1691.9999
$ws.Range("N103").Value = -2108.75
$ws.Range("H132").Value = 4251.72
$ws.Range("I132").Value = 4649.727
$ws.Range("K132").Value = 13949.181
$ws.Range("M132").Value = -11419.181
$ws.Range("H137").Value = 62209.6
$ws.Range("I137").Value = 101198.22
$ws.Range("J137").Value = 3726.6667
$ws.Range("K137").Value = 303594.66
$ws.Range("L137").Value = 11180.0001
$ws.Range("M137").Value = -301044.66
$ws.Range("N137").Value = -16280.0001
$ws.Range("H138").Value = 3242.82
$ws.Range("J138").Value = 4104.484
$ws.Range("L138").Value = 12313.452
$ws.Range("N138").Value = -22593.452
$ws.Range("H139").Value = 95000
$ws.Range("J139").Value = 95000
$ws.Range("L139").Value = 95000
$ws.Range("N139").Value = -105280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3731.4546
$ws.Range("I102").Value = 3004.2942
$ws.Range("J102").Value = 6203.8
$ws.Range("K102").Value = 3004.2942
$ws.Range("L102").Value = 6203.8
$ws.Range("M102").Value = -1382.2942
$ws.Range("N102").Value = -9447.799999999999
$ws.Range("H122").Value = 8232886.5
$ws.Range("I122").Value = 10102952
$ws.Range("J122").Value = 4598.6
$ws.Range("K122").Value = 30308856
$ws.Range("L122").Value = 13795.8
$ws.Range("M122").Value = -30306406
$ws.Range("N122").Value = -18695.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 837
$ws.Range("J64").Value = 731.4
$ws.Range("L64").Value = 731.4
$ws.Range("N64").Value = -1181.4
$ws.Range("H67").Value = 837
$ws.Range("J67").Value = 731.4
$ws.Range("L67").Value = 731.4
$ws.Range("N67").Value = -2291.4
$ws.Range("H86").Value = 4274.718
$ws.Range("I86").Value = 6054.68
$ws.Range("J86").Value = 1096.2142
$ws.Range("K86").Value = 6054.68
$ws.Range("L86").Value = 1096.2142
$ws.Range("M86").Value = -4931.68
$ws.Range("N86").Value = -3342.2142
$ws.Range("H89").Value = 4274.718
$ws.Range("I89").Value = 6054.68
$ws.Range("J89").Value = 1096.2142
$ws.Range("K89").Value = 30273.4
$ws.Range("L89").Value = 5481.071
$ws.Range("M89").Value = -24657.4
$ws.Range("N89").Value = -16713.071
$ws.Range("H134").Value = 3959.1892
$ws.Range("I134").Value = 2102.4546
$ws.Range("K134").Value = 6307.3638
$ws.Range("M134").Value = -3772.3638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15024.279
$ws.Range("I31").Value = 1913.25
$ws.Range("J31").Value = 16369
$ws.Range("K31").Value = 1913.25
$ws.Range("L31").Value = 16369
$ws.Range("M31").Value = -1618.25
$ws.Range("N31").Value = -16959
$ws.Range("H34").Value = 15024.279
$ws.Range("I34").Value = 1913.25
$ws.Range("J34").Value = 16369
$ws.Range("K34").Value = 1913.25
$ws.Range("L34").Value = 16369
$ws.Range("M34").Value = -1711.25
$ws.Range("N34").Value = -16773
$ws.Range("H58").Value = 1715.2858
$ws.Range("I58").Value = 1409.5454
$ws.Range("J58").Value = 2836.3333
$ws.Range("K58").Value = 1409.5454
$ws.Range("L58").Value = 2836.3333
$ws.Range("M58").Value = -1206.5454
$ws.Range("N58").Value = -3242.3333
$ws.Range("H86").Value = 9499.297
$ws.Range("I86").Value = 7432.7144
$ws.Range("J86").Value = 12211.6875
$ws.Range("K86").Value = 7432.7144
$ws.Range("L86").Value = 12211.6875
$ws.Range("M86").Value = -6309.7144
$ws.Range("N86").Value = -14457.6875
$ws.Range("H89").Value = 9499.297
$ws.Range("I89").Value = 7432.7144
$ws.Range("J89").Value = 12211.6875
$ws.Range("K89").Value = 37163.572
$ws.Range("L89").Value = 61058.4375
$ws.Range("M89").Value = -31547.572
$ws.Range("N89").Value = -72290.4375
$ws.Range("H120").Value = 96665
$ws.Range("J120").Value = 96665
$ws.Range("L120").Value = 96665
$ws.Range("N120").Value = -103923
$ws.Range("H134").Value = 2929.0417
$ws.Range("I134").Value = 1943.0625
$ws.Range("K134").Value = 5829.1875
$ws.Range("M134").Value = -3294.1875
$ws.Range("H136").Value = 1715.2858
$ws.Range("I136").Value = 1409.5454
$ws.Range("J136").Value = 2836.3333
$ws.Range("K136").Value = 4228.6362
$ws.Range("L136").Value = 8508.999899999999
$ws.Range("M136").Value = -1678.6362
$ws.Range("N136").Value = -13608.9999
$ws.Range("H138").Value = 98900
$ws.Range("J138").Value = 157800
$ws.Range("L138").Value = 157800
$ws.Range("N138").Value = -168080

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 247.71428
$ws.Range("I23").Value = 38
$ws.Range("K23").Value = 114
$ws.Range("M23").Value = 121
$ws.Range("H92").Value = 210
$ws.Range("J92").Value = 210
$ws.Range("L92").Value = 630
$ws.Range("N92").Value = -3126
$ws.Range("H134").Value = 7833
$ws.Range("I134").Value = 3750
$ws.Range("K134").Value = 11250
$ws.Range("M134").Value = -6180

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1001.7
$ws.Range("I97").Value = 1150.6875
$ws.Range("J97").Value = 405.75
$ws.Range("K97").Value = 1150.6875
$ws.Range("L97").Value = 405.75
$ws.Range("M97").Value = -654.6875
$ws.Range("N97").Value = -1397.75
$ws.Range("H128").Value = 93998
$ws.Range("J128").Value = 93998
$ws.Range("L128").Value = 93998
$ws.Range("N128").Value = -103958
$ws.Range("H132").Value = 3147
$ws.Range("I132").Value = 2474.6072
$ws.Range("K132").Value = 7423.821599999999
$ws.Range("M132").Value = -4893.821599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 17332.834
$ws.Range("I48").Value = 8499.25
$ws.Range("K48").Value = 8499.25
$ws.Range("M48").Value = -7838.25
$ws.Range("H68").Value = 2131.182
$ws.Range("I68").Value = 1465.8889
$ws.Range("K68").Value = 1465.8889
$ws.Range("M68").Value = -716.8888999999999
$ws.Range("H71").Value = 2131.182
$ws.Range("I71").Value = 1465.8889
$ws.Range("K71").Value = 7329.4445
$ws.Range("M71").Value = -3585.4445
$ws.Range("H82").Value = 1582
$ws.Range("I82").Value = 1181.3334
$ws.Range("J82").Value = 1782.3334
$ws.Range("K82").Value = 1181.3334
$ws.Range("L82").Value = 1782.3334
$ws.Range("M82").Value = -820.3334
$ws.Range("N82").Value = -2504.3334
$ws.Range("H85").Value = 1582
$ws.Range("I85").Value = 1181.3334
$ws.Range("J85").Value = 1782.3334
$ws.Range("K85").Value = 1181.3334
$ws.Range("L85").Value = 1782.3334
$ws.Range("M85").Value = 66.66660000000002
$ws.Range("N85").Value = -4278.3334
$ws.Range("H122").Value = 8096.625
$ws.Range("I122").Value = 4497.5
$ws.Range("K122").Value = 13492.5
$ws.Range("M122").Value = -11042.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 17060.715
$ws.Range("J54").Value = 18000
$ws.Range("L54").Value = 18000
$ws.Range("N54").Value = -19040
$ws.Range("H80").Value = 50000
$ws.Range("I80").Value = 50000
$ws.Range("K80").Value = 50000
$ws.Range("M80").Value = -49002
$ws.Range("H83").Value = 50000
$ws.Range("I83").Value = 50000
$ws.Range("K83").Value = 150000
$ws.Range("M83").Value = -145008
$ws.Range("H132").Value = 48310.41
$ws.Range("I132").Value = 9934.134
$ws.Range("J132").Value = 130545.29
$ws.Range("K132").Value = 29802.402
$ws.Range("L132").Value = 391635.87
$ws.Range("M132").Value = -27272.402
$ws.Range("N132").Value = -396695.87
$ws.Range("H136").Value = 3090.1333
$ws.Range("I136").Value = 2441.0908
$ws.Range("J136").Value = 4875
$ws.Range("K136").Value = 7323.2724
$ws.Range("L136").Value = 14625
$ws.Range("M136").Value = -4773.2724
$ws.Range("N136").Value = -19725
